$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 371
$ws.Cells.Item(371, 6).Value2 = 'Edmonton Oil Kings'
$ws.Cells.Item(371, 7).Value2 = 0

# Row 372
$ws.Cells.Item(372, 6).Value2 = 'Saskatoon Blades'
$ws.Cells.Item(372, 7).Value2 = 1

# Row 373
$ws.Cells.Item(373, 5).Value2 = 'Portland Winterhawks'
$ws.Cells.Item(373, 6).Value2 = 'Portland Winterhawks'
$ws.Cells.Item(373, 7).Value2 = 1

# Row 374
$ws.Cells.Item(374, 6).Value2 = 'Seattle Thunderbirds'
$ws.Cells.Item(374, 7).Value2 = 1

# Row 375
$ws.Cells.Item(375, 6).Value2 = 'Everett Silvertips'
$ws.Cells.Item(375, 7).Value2 = 0

# Row 376
$ws.Cells.Item(376, 1).Value2 = 1021923
$ws.Cells.Item(376, 2).Value2 = 'Tue, Mar 18, 2025'
$ws.Cells.Item(376, 3).Value2 = 'Prince Albert Raiders'
$ws.Cells.Item(376, 4).Value2 = 'Moose Jaw Warriors'
$ws.Cells.Item(376, 5).Value2 = 'Prince Albert Raiders'
$ws.Cells.Item(376, 6).Value2 = 'Prince Albert Raiders'
$ws.Cells.Item(376, 7).Value2 = 1

# Row 377
$ws.Cells.Item(377, 1).Value2 = 1021924
$ws.Cells.Item(377, 2).Value2 = 'Tue, Mar 18, 2025'
$ws.Cells.Item(377, 3).Value2 = 'Regina Pats'
$ws.Cells.Item(377, 4).Value2 = 'Swift Current Broncos'
$ws.Cells.Item(377, 5).Value2 = 'Swift Current Broncos'
$ws.Cells.Item(377, 6).Value2 = 'Regina Pats'
$ws.Cells.Item(377, 7).Value2 = 0

# Row 378
$ws.Cells.Item(378, 1).Value2 = 1021922
$ws.Cells.Item(378, 2).Value2 = 'Tue, Mar 18, 2025'
$ws.Cells.Item(378, 3).Value2 = 'Kamloops Blazers'
$ws.Cells.Item(378, 4).Value2 = 'Victoria Royals'
$ws.Cells.Item(378, 5).Value2 = 'Victoria Royals'
$ws.Cells.Item(378, 6).Value2 = 'Victoria Royals'
$ws.Cells.Item(378, 7).Value2 = 1

# Row 379
$ws.Cells.Item(379, 1).Value2 = 1021925
$ws.Cells.Item(379, 2).Value2 = 'Wed, Mar 19, 2025'
$ws.Cells.Item(379, 3).Value2 = 'Brandon Wheat Kings'
$ws.Cells.Item(379, 4).Value2 = 'Saskatoon Blades'
$ws.Cells.Item(379, 5).Value2 = 'Saskatoon Blades'
$ws.Cells.Item(379, 6).Value2 = 'Saskatoon Blades'
$ws.Cells.Item(379, 7).Value2 = 1

# Row 380
$ws.Cells.Item(380, 1).Value2 = 1021926
$ws.Cells.Item(380, 2).Value2 = 'Wed, Mar 19, 2025'
$ws.Cells.Item(380, 3).Value2 = 'Calgary Hitmen'
$ws.Cells.Item(380, 4).Value2 = 'Red Deer Rebels'
$ws.Cells.Item(380, 5).Value2 = 'Calgary Hitmen'
$ws.Cells.Item(380, 6).Value2 = 'Calgary Hitmen'
$ws.Cells.Item(380, 7).Value2 = 1

# Row 381
$ws.Cells.Item(381, 1).Value2 = 1021927
$ws.Cells.Item(381, 2).Value2 = 'Wed, Mar 19, 2025'
$ws.Cells.Item(381, 3).Value2 = 'Kelowna Rockets'
$ws.Cells.Item(381, 4).Value2 = 'Victoria Royals'
$ws.Cells.Item(381, 5).Value2 = 'Victoria Royals'
$ws.Cells.Item(381, 6).Value2 = 'Victoria Royals'
$ws.Cells.Item(381, 7).Value2 = 1

# Row 382
$ws.Cells.Item(382, 1).Value2 = 1021928
$ws.Cells.Item(382, 2).Value2 = 'Fri, Mar 21, 2025'
$ws.Cells.Item(382, 3).Value2 = 'Brandon Wheat Kings'
$ws.Cells.Item(382, 4).Value2 = 'Regina Pats'
$ws.Cells.Item(382, 5).Value2 = 'Regina Pats'

# Row 383
$ws.Cells.Item(383, 1).Value2 = 1021931
$ws.Cells.Item(383, 2).Value2 = 'Fri, Mar 21, 2025'
$ws.Cells.Item(383, 3).Value2 = 'Lethbridge Hurricanes'
$ws.Cells.Item(383, 4).Value2 = 'Edmonton Oil Kings'
$ws.Cells.Item(383, 5).Value2 = 'Edmonton Oil Kings'

# Row 384
$ws.Cells.Item(384, 1).Value2 = 1021932
$ws.Cells.Item(384, 2).Value2 = 'Fri, Mar 21, 2025'
$ws.Cells.Item(384, 3).Value2 = 'Swift Current Broncos'
$ws.Cells.Item(384, 4).Value2 = 'Moose Jaw Warriors'
$ws.Cells.Item(384, 5).Value2 = 'Moose Jaw Warriors'

# Row 385
$ws.Cells.Item(385, 1).Value2 = 1021935
$ws.Cells.Item(385, 2).Value2 = 'Fri, Mar 21, 2025'
$ws.Cells.Item(385, 3).Value2 = 'Saskatoon Blades'
$ws.Cells.Item(385, 4).Value2 = 'Prince Albert Raiders'
$ws.Cells.Item(385, 5).Value2 = 'Saskatoon Blades'

# Row 386
$ws.Cells.Item(386, 1).Value2 = 1021930
$ws.Cells.Item(386, 2).Value2 = 'Fri, Mar 21, 2025'
$ws.Cells.Item(386, 3).Value2 = 'Kamloops Blazers'
$ws.Cells.Item(386, 4).Value2 = 'Vancouver Giants'
$ws.Cells.Item(386, 5).Value2 = 'Kamloops Blazers'

# Row 387
$ws.Cells.Item(387, 1).Value2 = 1021933
$ws.Cells.Item(387, 2).Value2 = 'Fri, Mar 21, 2025'
$ws.Cells.Item(387, 3).Value2 = 'Portland Winterhawks'
$ws.Cells.Item(387, 4).Value2 = 'Tri-City Americans'
$ws.Cells.Item(387, 5).Value2 = 'Portland Winterhawks'

# Row 388
$ws.Cells.Item(388, 1).Value2 = 1021934
$ws.Cells.Item(388, 2).Value2 = 'Fri, Mar 21, 2025'
$ws.Cells.Item(388, 3).Value2 = 'Prince George Cougars'
$ws.Cells.Item(388, 4).Value2 = 'Victoria Royals'
$ws.Cells.Item(388, 5).Value2 = 'Victoria Royals'

# Row 389
$ws.Cells.Item(389, 1).Value2 = 1021929
$ws.Cells.Item(389, 2).Value2 = 'Fri, Mar 21, 2025'
$ws.Cells.Item(389, 3).Value2 = 'Everett Silvertips'
$ws.Cells.Item(389, 4).Value2 = 'Wenatchee Wild'
$ws.Cells.Item(389, 5).Value2 = 'Everett Silvertips'

# Row 390
$ws.Cells.Item(390, 1).Value2 = 1021936
$ws.Cells.Item(390, 2).Value2 = 'Fri, Mar 21, 2025'
$ws.Cells.Item(390, 3).Value2 = 'Spokane Chiefs'
$ws.Cells.Item(390, 4).Value2 = 'Seattle Thunderbirds'
$ws.Cells.Item(390, 5).Value2 = 'Spokane Chiefs'

# Update selection / view position to match the post-edit Excel state
$win = $excel.ActiveWindow
$win.ScrollRow = 364
$null = $ws.Range("H383").Select()
